# Add two new columns: I ("I0") and J ("IF")
# I0 is always 1, IF duplicates the value already present in column H (IP)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used data row based on column A (rows 2..22 hold data)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# --- Headers (row 1) ---
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Copy the formatting (style) of the existing H1 header cell onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows ---
for ($r = 2; $r -le $lastRow; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
